# Updates the "cryptos" listing (Sheet1) with refreshed price/volume data.
# Note: price values that look like plain numbers (e.g. "96.72") are written
# with a leading apostrophe to force Excel to keep them as text, matching
# the original inline-string cell content (these are decimal-formatted
# price strings, not numeric values - some even use "." as a thousands
# separator, e.g. "44.045.02", which would not round-trip as a number).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '44.045.02'
$ws.Range('E2').Value = '  +1.28%  '
$ws.Range('D3').Value = '2.254.12'
$ws.Range('E3').Value = '  +2.85%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('B5').Value = 'Solana'
$ws.Range('C5').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D5').Value = "'96.72"
$ws.Range('E5').Value = '  +17.10%  '
$ws.Range('B6').Value = 'BNB'
$ws.Range('C6').Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range('D6').Value = "'271.79"
$ws.Range('E6').Value = '  +5.96%  '
$ws.Range('D7').Value = "'0.624"
$ws.Range('E7').Value = '  +1.01%  '
$ws.Range('E8').Value = '  -0.08%  '
$ws.Range('D9').Value = "'0.654"
$ws.Range('E9').Value = '  +10.91%  '
$ws.Range('D10').Value = "'47.67"
$ws.Range('E10').Value = '  +7.71%  '
$ws.Range('D11').Value = "'0.0957"
$ws.Range('E11').Value = '  +4.45%  '
$ws.Range('D12').Value = "'8.48"
$ws.Range('E12').Value = '  +18.69%  '
$ws.Range('E13').Value = '  +1.35%  '
$ws.Range('D14').Value = "'15.47"
$ws.Range('E14').Value = '  +7.97%  '
$ws.Range('D15').Value = '2.583.91'
$ws.Range('E15').Value = '  +2.67%  '
$ws.Range('D16').Value = "'0.834"
$ws.Range('E16').Value = '  +7.41%  '
$ws.Range('D17').Value = '2.250.90'
$ws.Range('E17').Value = '  +0.63%  '
$ws.Range('D18').Value = '44.012.45'
$ws.Range('E18').Value = '  +1.39%  '
$ws.Range('D19').Value = "'0.0000107"
$ws.Range('E19').Value = '  +4.13%  '
$ws.Range('D20').Value = "'6.27"
$ws.Range('E20').Value = '  +6.31%  '
$ws.Range('D21').Value = "'71.29"
$ws.Range('E21').Value = '  +3.21%  '
$ws.Range('D22').Value = "'2.34"
$ws.Range('E22').Value = '  -1.46%  '
$ws.Range('D23').Value = "'236.34"
$ws.Range('E23').Value = '  +2.38%  '
$ws.Range('D24').Value = "'9.26"
$ws.Range('E24').Value = '  +7.00%  '
$ws.Range('B25').Value = 'Cosmos'
$ws.Range('C25').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D25').Value = "'11.68"
$ws.Range('E25').Value = '  +10.17%  '
$ws.Range('B26').Value = 'Dai'
$ws.Range('C26').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D26').Value = "'1.00"
$ws.Range('E26').Value = '  -0.01%  '
$ws.Range('D27').Value = "'2.52"
$ws.Range('E27').Value = '  +13.90%  '
$ws.Range('D28').Value = "'3.50"
$ws.Range('E28').Value = '  +2.73%  '
$ws.Range('D29').Value = "'40.96"
$ws.Range('E29').Value = '  +1.09%  '
$ws.Range('D30').Value = "'2.26"
$ws.Range('E30').Value = '  +0.06%  '
$ws.Range('D31').Value = "'172.40"
$ws.Range('E31').Value = '  -0.69%  '
$ws.Range('D32').Value = "'0.0929"
$ws.Range('E32').Value = '  +7.39%  '
$ws.Range('D33').Value = "'21.13"
$ws.Range('E33').Value = '  +3.92%  '
$ws.Range('D34').Value = "'5.60"
$ws.Range('E34').Value = '  +5.19%  '
$ws.Range('B35').Value = 'Kaspa'
$ws.Range('C35').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D35').Value = "'0.115"
$ws.Range('E35').Value = '  +2.04%  '
$ws.Range('B36').Value = 'Stellar'
$ws.Range('C36').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D36').Value = "'0.125"
$ws.Range('E36').Value = '  +1.92%  '
$ws.Range('D37').Value = "'0.0357"
$ws.Range('E37').Value = '  -0.88%  '
$ws.Range('D38').Value = "'4.41"
$ws.Range('E38').Value = '  -0.72%  '
$ws.Range('D39').Value = "'3.70"
$ws.Range('E39').Value = '  +30.48%  '
$ws.Range('D40').Value = "'0.234"
$ws.Range('E40').Value = '  +18.15%  '
$ws.Range('D41').Value = "'13.23"
$ws.Range('E41').Value = '  +7.05%  '
$ws.Range('D42').Value = "'2.19"
$ws.Range('E42').Value = '  +4.77%  '
$ws.Range('D43').Value = "'62.96"
$ws.Range('E43').Value = '  +0.07%  '
$ws.Range('D44').Value = "'5.47"
$ws.Range('E44').Value = '  +0.57%  '
$ws.Range('B45').Value = 'Cronos'
$ws.Range('C45').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D45').Value = "'0.101"
$ws.Range('E45').Value = '  +4.03%  '
$ws.Range('B46').Value = 'Aave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D46').Value = "'102.36"
$ws.Range('E46').Value = '  +2.24%  '
$ws.Range('B47').Value = 'FraxShare'
$ws.Range('C47').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D47').Value = "'8.49"
$ws.Range('E47').Value = '  +3.45%  '
$ws.Range('D48').Value = "'1.18"
$ws.Range('E48').Value = '  +7.09%  '
$ws.Range('E49').Value = '  +2.22%  '
$ws.Range('D50').Value = "'0.450"
$ws.Range('E50').Value = '  +2.00%  '
$ws.Range('D51').Value = '2.469.42'
$ws.Range('E51').Value = '  +2.52%  '
